# decayCHAIN.xlsx — insert a new "Reference particle" row for a second
# species (muon) right above the existing "Vacuum chamber" row, mirroring
# the existing "Species" / "pion" row (row 4) but as "Species 1" / "muon".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("decayCHAIN")

# Insert a new row above row 6 (current "Vacuum chamber" row), shifting
# everything below it down by one.
$ws.Rows.Item(6).Insert()

# Copy the formatting of row 4 (the existing "Species"/"pion" row) onto the
# freshly inserted row 6 (restricted to the used columns, A:H), then
# overwrite the text with the new species.
$ws.Range("A4:H4").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "Facility"
$ws.Range("C6").Value = "Global"
$ws.Range("D6").Value = "Reference particle"
$ws.Range("E6").Value = "Species 1"
$ws.Range("F6").Value = "muon"
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null

$ws.Range("F6").Select()
